# Add two new "fade in" click-triggered animation groups to slide 2
# (Machine Learning - PCA.pptx):
#   1) Shape id=8  (clickEffect, fade)  together-with  Shape id=4  (withEffect, fade)
#   2) Shape id=17 (clickEffect, fade)  together-with  Shape id=18 (withEffect, fade)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$tl = $s.TimeLine
$ms = $tl.MainSequence

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

$shape8  = Get-ShapeById $s 8
$shape4  = Get-ShapeById $s 4
$shape17 = Get-ShapeById $s 17
$shape18 = Get-ShapeById $s 18

# msoAnimEffectFade = 10, msoAnimTriggerOnPageClick = 1, msoAnimTriggerWithPrevious = 2

# First group: click on shape 8, shape 4 appears with it
$eff1 = $ms.AddEffect($shape8, 10, 0, 1)
$eff2 = $ms.AddEffect($shape4, 10, 0, 2)

# Second group: click on shape 17, shape 18 appears with it
$eff3 = $ms.AddEffect($shape17, 10, 0, 1)
$eff4 = $ms.AddEffect($shape18, 10, 0, 2)
